$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update label to clarify the unit of the wheel diameter
$ws.Range("D12").Value = "Diameter (inches)"

# New wheel tread / diameter measurement
$ws.Range("E12").Value = 5

# Move selection to match the saved state (E13)
$ws.Range("E13").Select()

$wb.Save()
